{"js": "// Replace each three-digit-by-one-digit multiplication expression with its\n// updated counterpart. Every occurrence is unique in the document, so a\n// simple search-and-replace per pair is safe and order-independent.\nconst replacements = [\n  [\"673\u00d79=\", \"422\u00d73=\"],\n  [\"976\u00d72=\", \"268\u00d75=\"],\n  [\"379\u00d75=\", \"644\u00d79=\"],\n  [\"909\u00d79=\", \"331\u00d79=\"],\n  [\"565\u00d73=\", \"191\u00d73=\"],\n  [\"758\u00d72=\", \"325\u00d73=\"],\n  [\"813\u00d76=\", \"608\u00d74=\"],\n  [\"714\u00d76=\", \"870\u00d72=\"],\n  [\"557\u00d73=\", \"649\u00d72=\"],\n  [\"357\u00d74=\", \"624\u00d73=\"],\n  [\"214\u00d76=\", \"243\u00d78=\"],\n  [\"502\u00d76=\", \"841\u00d75=\"],\n  [\"113\u00d74=\", \"609\u00d76=\"],\n  [\"288\u00d74=\", \"886\u00d79=\"],\n  [\"793\u00d75=\", \"831\u00d76=\"],\n  [\"138\u00d72=\", \"210\u00d72=\"],\n  [\"275\u00d79=\", \"889\u00d75=\"],\n  [\"448\u00d77=\", \"786\u00d74=\"],\n  [\"404\u00d78=\", \"179\u00d77=\"],\n  [\"267\u00d74=\", \"663\u00d73=\"],\n  [\"125\u00d79=\", \"166\u00d76=\"],\n  [\"769\u00d78=\", \"599\u00d78=\"],\n  [\"853\u00d78=\", \"300\u00d76=\"],\n  [\"913\u00d78=\", \"520\u00d74=\"],\n  [\"575\u00d73=\", \"453\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression with its\n# updated counterpart. Every \"old\" string occurs exactly once in the\n# document, so Find/Replace with ReplaceAll is safe per pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"673\u00d79=\"; New = \"422\u00d73=\"},\n    @{Old = \"976\u00d72=\"; New = \"268\u00d75=\"},\n    @{Old = \"379\u00d75=\"; New = \"644\u00d79=\"},\n    @{Old = \"909\u00d79=\"; New = \"331\u00d79=\"},\n    @{Old = \"565\u00d73=\"; New = \"191\u00d73=\"},\n    @{Old = \"758\u00d72=\"; New = \"325\u00d73=\"},\n    @{Old = \"813\u00d76=\"; New = \"608\u00d74=\"},\n    @{Old = \"714\u00d76=\"; New = \"870\u00d72=\"},\n    @{Old = \"557\u00d73=\"; New = \"649\u00d72=\"},\n    @{Old = \"357\u00d74=\"; New = \"624\u00d73=\"},\n    @{Old = \"214\u00d76=\"; New = \"243\u00d78=\"},\n    @{Old = \"502\u00d76=\"; New = \"841\u00d75=\"},\n    @{Old = \"113\u00d74=\"; New = \"609\u00d76=\"},\n    @{Old = \"288\u00d74=\"; New = \"886\u00d79=\"},\n    @{Old = \"793\u00d75=\"; New = \"831\u00d76=\"},\n    @{Old = \"138\u00d72=\"; New = \"210\u00d72=\"},\n    @{Old = \"275\u00d79=\"; New = \"889\u00d75=\"},\n    @{Old = \"448\u00d77=\"; New = \"786\u00d74=\"},\n    @{Old = \"404\u00d78=\"; New = \"179\u00d77=\"},\n    @{Old = \"267\u00d74=\"; New = \"663\u00d73=\"},\n    @{Old = \"125\u00d79=\"; New = \"166\u00d76=\"},\n    @{Old = \"769\u00d78=\"; New = \"599\u00d78=\"},\n    @{Old = \"853\u00d78=\"; New = \"300\u00d76=\"},\n    @{Old = \"913\u00d78=\"; New = \"520\u00d74=\"},\n    @{Old = \"575\u00d73=\"; New = \"453\u00d75=\"}\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $pair.New, 2) | Out-Null\n}\n"}
